$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.741.17"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.848.29"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4335"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3653"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07336"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8764"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.73"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").Value = "1.802.16"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.346"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.519"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06934"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008994"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("D22").Value = "27.607.48"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.977"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.12%  "

$ws.Range("D25").Value = "2.050.64"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.984"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "120.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +8.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.251"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.861"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08913"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7529"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.535"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.962"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.121"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.109"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05410"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01931"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.836"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5085"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.652"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.319"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06544"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4659"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.37"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.621"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.07%  "
